$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/add header row: Ano | Cadastrado | Sem Cadastro
$ws.Range("B1").Value = "Cadastrado"

# Copy the existing header formatting (bold, centered, bordered) onto the new column
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Sem Cadastro"

# Data: Ano, Cadastrado, Sem Cadastro
$data = @(
    @(2020, 17536.47, 6329.94),
    @(2021, 252101.08, 27817.06),
    @(2022, 489664.5, 14684.78),
    @(2023, 791532.51, 10932.38),
    @(2024, 1184081.5, 12980.97),
    @(2025, 391120.88, 5486.44)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
